$d = $word.ActiveDocument

# Locate the paragraph that ends with "... più credibile" (the last bullet
# under "Alla verifica:") so we can add a new sibling bullet right after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*rendere l'applicazione più credibile*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*credibile*") {
            $target = $p
            break
        }
    }
}

# Insert a brand-new paragraph right after the target one; it automatically
# inherits the same pPr (Paragrafoelenco style, ilvl=1, numId=3) as Word
# does when you press Enter at the end of a list item.
$target.Range.InsertParagraphAfter()
$firstNew = $target.Next()
$firstNew.Range.InsertAfter("Dentro il codice ho inserito tutti gli appunti logici/per tesi/di riferimento d")

# Add the second run's text as its own paragraph first (so it becomes an
# independent run rather than being coalesced into the previous one), then
# fold it back into the first paragraph by deleting the paragraph mark that
# separates them. The result is a single paragraph made of two runs, just
# like the target revision.
$firstNew.Range.InsertParagraphAfter()
$secondNew = $firstNew.Next()
$secondNew.Range.InsertAfter("i tutto l’ambaradan realizzato")

$mark = $d.Range($firstNew.Range.End - 1, $firstNew.Range.End)
$mark.Delete()
